# This workbook originally carried the indicator metadata in three columns:
#   A = English label, B = Russian label, C = Russian content.
# The edit drops the English-language column A entirely; the Russian
# label/content columns shift left to become the new A/B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stash a copy of the hyperlinked cell's formatting (it lives at C8, which
# becomes B8 once column A is removed) in a scratch cell far outside the
# used range, so we can restore it after re-pointing the hyperlink below.
$ws.Range("C8").Copy($ws.Range("Z1"))

# Delete column A (English labels); B->A and C->B shift left automatically,
# and shared strings / column widths / row data all move with them.
$ws.Columns("A").Delete()

# The workbook's lone hyperlink (the contact e-mail) was anchored on the
# old C8. Column deletion does not re-anchor it, so re-create it on the
# cell it now occupies, B8.
$ws.Range("C8").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B8"), "mailto:Jmaatkulova@stat.kg")

# Adding the hyperlink reset B8's formatting to the default hyperlink
# style; restore the original look (border/wrap/etc.) from the stashed
# copy. The scratch cell shifted left along with everything else when
# column A was deleted, so it is now Y1, not Z1.
$ws.Range("Y1").Copy($ws.Range("B8"))
$ws.Range("Y1").Clear()

# Match the saved selection/active cell of the edited workbook.
$ws.Range("E3").Select()
